$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Text clean-up / typo fixes (merges adjacent runs that had identical
#    formatting back into a single run, and fixes a few "Herrn" -> "Herr"
#    typos plus the "insterlliert" -> "installiert" misspelling).
# ---------------------------------------------------------------------------

# "Ziel des Angriffs: Lahm legen einer Webseite, " + "Server etc." -> one run
$d.Content.Find.Execute(": Lahm legen einer Webseite, Server etc.", $false, $false, $false, $false, $false, $true, 1, $false, ": Lahm legen einer Webseite, Server etc.", 2) | Out-Null

# "Akteure : Herr Mustermann, " + "Herr Stern" -> one run
$d.Content.Find.Execute(" : Herr Mustermann, Herr Stern", $false, $false, $false, $false, $false, $true, 1, $false, " : Herr Mustermann, Herr Stern", 2) | Out-Null

# Fallbeschreibung paragraph: merge many runs into one + fix "Herrn Stern" -> "Herr Stern"
$d.Content.Find.Execute("Herr Mustermann möchte den Server von Herrn Stern lahm legen, damit Herr Stern nicht arbeiten kann. Weil Herr Stern den letzten Keks genommen hat. ", $false, $false, $false, $false, $false, $true, 1, $false, "Herr Mustermann möchte den Server von Herr Stern lahm legen, damit Herr Stern nicht arbeiten kann. Weil Herr Stern den letzten Keks genommen hat. ", 2) | Out-Null

# "2. " + "... von Herrn Stern ..." -> one run + fix "Herrn Stern" -> "Herr Stern"
$d.Content.Find.Execute("2. Herr Mustermann kennt die IP des Servers von Herrn Stern und trägt diese ein :", $false, $false, $false, $false, $false, $true, 1, $false, "2. Herr Mustermann kennt die IP des Servers von Herr Stern und trägt diese ein :", 2) | Out-Null

# "3. " + "Angriff starten:" -> one run
$d.Content.Find.Execute("3. Angriff starten:", $false, $false, $false, $false, $false, $true, 1, $false, "3. Angriff starten:", 2) | Out-Null

# "4. " + "Jetzt wird der Server von Herrn Stern mit Anfragen ..." -> one run (no typo fix here)
$d.Content.Find.Execute("4. Jetzt wird der Server von Herrn Stern mit Anfragen", $false, $false, $false, $false, $false, $true, 1, $false, "4. Jetzt wird der Server von Herrn Stern mit Anfragen", 2) | Out-Null

# Final paragraph: merge the trailing "." run in + fix "Herrn Stern" -> "Herr Stern"
$d.Content.Find.Execute("Der Server von Herrn Stern läuft nun langsamer da er mit Anfragen „bombardiert“ wird und dadurch wird das Arbeiten mit dem Server erschwert.", $false, $false, $false, $false, $false, $true, 1, $false, "Der Server von Herr Stern läuft nun langsamer da er mit Anfragen „bombardiert“ wird und dadurch wird das Arbeiten mit dem Server erschwert.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. "insterlliert" -> "installiert" : fix the typo by replacing the "erl"
#    chunk with "al", and keep the inserted "al" as its own run (matching
#    the authored diff's run split) by nudging a character property on it.
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("insterlliert", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$s = $rng.Start
$e = $rng.End

$mid = $d.Range($s + 4, $s + 7)
$mid.Text = "al"

$alRng = $d.Range($s + 4, $s + 6)
$alRng.Font.Bold = $true
$alRng.Font.Bold = $false

# ---------------------------------------------------------------------------
# 3. Table indentation / cell margin tweaks (tblInd, tblCellMar, tcMar).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $t = $d.Tables.Item($i)

    if ($t.Rows.LeftIndent -eq -0.25) {
        $t.Rows.LeftIndent = -0.5
    } elseif ($t.Rows.LeftIndent -eq 1.65) {
        $t.Rows.LeftIndent = 1.4
    }

    $t.LeftPadding = 4.9
    foreach ($c in $t.Range.Cells) {
        $c.LeftPadding = 4.9
    }
}

Write-Output "done"
